$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.282.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.910.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.35"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4723"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.97%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.74"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08028"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.50"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.924.73"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.882"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.121"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.56"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06633"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.62"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.307.47"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.523"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.198"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.179.68"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.73"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.987"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +10.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.105"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.41"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09492"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.423"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.539"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.379"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06075"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02247"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.226"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.171"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5862"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.523"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +11.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1834"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.09"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07862"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.273"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5522"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.08"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.06"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.22"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.50%  "
